$d = $word.ActiveDocument

# --- Step 1: Row 8, cell 2 - remove the _GoBack bookmark that currently
#     sits at the end of "Organizar la galería de productos (tabla y
#     tamaño de imágenes)" (keep the text exactly as-is).
$tbl = $d.Tables.Item(1)
$cell8 = $tbl.Cell(8, 2)
$rng8 = $cell8.Range
$rng8.Collapse(0)

$xmlRow8 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p w:rsidR="00906F2F" w:rsidRPr="00906F2F" w:rsidRDefault="00906F2F" w:rsidP="007D6730">
<w:r w:rsidRPr="00906F2F"><w:t>Organizar la galería de product</w:t></w:r>
<w:r w:rsidR="00592D8F"><w:t>os (tabla y tamaño de imágenes)</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng8.InsertXML($xmlRow8)

# --- Step 2: Row 3, cell 2 - after "inscribir cursos", append a space,
#     the (now relocated) _GoBack bookmark, and the cyan-highlighted
#     "(falta poner selcts)" note (with "selcts" flagged as a misspelling).
$tbl = $d.Tables.Item(1)
$cell3 = $tbl.Cell(3, 2)
$rng3 = $cell3.Range
$found = $rng3.Find.Execute("inscribir cursos", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

$xmlRow3 = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:r><w:t xml:space="preserve"> </w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
<w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t xml:space="preserve">(falta poner </w:t></w:r>
<w:proofErr w:type="spellStart"/>
<w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>selcts</w:t></w:r>
<w:proofErr w:type="spellEnd"/>
<w:r><w:rPr><w:highlight w:val="cyan"/></w:rPr><w:t>)</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@
$rng3.InsertXML($xmlRow3)
